$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel;
# force them to stay Text (matching the original inlineStr cell type) by briefly
# switching to a Text number format for the assignment, then restoring the style
# so no visible formatting change is left behind.
$ws.Range('D2').Value = '62.984.41'
$ws.Range('E2').Value = '  +6.32%  '
$ws.Range('D3').Value = '3.112.83'
$ws.Range('E3').Value = '  +3.97%  '
$ws.Range('E4').Value = '  +0.03%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = "@"
$cell.Value = '587.28'
$cell.Style = "Normal"
$ws.Range('E5').Value = '  +4.84%  '
$cell = $ws.Range('D6')
$cell.NumberFormat = "@"
$cell.Value = '144.03'
$cell.Style = "Normal"
$ws.Range('E6').Value = '  +4.36%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = '3.101.79'
$ws.Range('E8').Value = '  +4.02%  '
$ws.Range('E9').Value = '  +2.56%  '
$ws.Range('E10').Value = '  +10.50%  '
$ws.Range('E11').Value = '  +10.45%  '
$ws.Range('E12').Value = '  +2.13%  '
$ws.Range('E13').Value = '  +6.75%  '
$cell = $ws.Range('D14')
$cell.NumberFormat = "@"
$cell.Value = '35.64'
$cell.Style = "Normal"
$ws.Range('E14').Value = '  +6.04%  '
$ws.Range('E15').Value = '  +0.50%  '
$ws.Range('D16').Value = '3.626.35'
$ws.Range('E16').Value = '  +4.01%  '
$ws.Range('E17').Value = '  -0.40%  '
$ws.Range('D18').Value = '62.941.21'
$ws.Range('E18').Value = '  +6.24%  '
$ws.Range('D19').Value = '3.107.20'
$ws.Range('E19').Value = '  +3.93%  '
$cell = $ws.Range('D20')
$cell.NumberFormat = "@"
$cell.Value = '453.88'
$cell.Style = "Normal"
$ws.Range('E20').Value = '  +5.79%  '
$cell = $ws.Range('D21')
$cell.NumberFormat = "@"
$cell.Value = '14.12'
$cell.Style = "Normal"
$ws.Range('E22').Value = '  +1.56%  '
$ws.Range('E23').Value = '  +6.80%  '
$cell = $ws.Range('D24')
$cell.NumberFormat = "@"
$cell.Value = '13.62'
$cell.Style = "Normal"
$ws.Range('E24').Value = '  +0.78%  '
$cell = $ws.Range('D25')
$cell.NumberFormat = "@"
$cell.Value = '82.01'
$cell.Style = "Normal"
$ws.Range('E25').Value = '  +2.14%  '
$ws.Range('E26').Value = '  +0.13%  '
$ws.Range('E27').Value = '  +2.11%  '
$ws.Range('E29').Value = '  +0.02%  '
$cell = $ws.Range('D30')
$cell.NumberFormat = "@"
$cell.Value = '8.24'
$cell.Style = "Normal"
$ws.Range('E30').Value = '  +5.27%  '
$cell = $ws.Range('D31')
$cell.NumberFormat = "@"
$cell.Value = '6.87'
$cell.Style = "Normal"
$ws.Range('E31').Value = '  +12.84%  '
$cell = $ws.Range('D33')
$cell.NumberFormat = "@"
$cell.Value = '26.99'
$cell.Style = "Normal"
$ws.Range('E33').Value = '  +5.02%  '
$ws.Range('E34').Value = '  +13.92%  '
$ws.Range('D35').Value = '0.0₃0809'
$ws.Range('E35').Value = '  +7.04%  '
$ws.Range('E36').Value = '  +4.18%  '
$cell = $ws.Range('D37')
$cell.NumberFormat = "@"
$cell.Value = '6.05'
$cell.Style = "Normal"
$ws.Range('E37').Value = '  +1.84%  '
$ws.Range('E38').Value = '  +13.03%  '
$cell = $ws.Range('D39')
$cell.NumberFormat = "@"
$cell.Value = '51.09'
$cell.Style = "Normal"
$ws.Range('E39').Value = '  +4.54%  '
$cell = $ws.Range('D40')
$cell.NumberFormat = "@"
$cell.Value = '8.76'
$cell.Style = "Normal"
$ws.Range('E40').Value = '  +1.15%  '
$cell = $ws.Range('D41')
$cell.NumberFormat = "@"
$cell.Value = '426.18'
$cell.Style = "Normal"
$ws.Range('E41').Value = '  +5.28%  '
$ws.Range('D42').Value = '2.970.80'
$ws.Range('E42').Value = '  +7.20%  '
$ws.Range('E43').Value = '  +5.76%  '
$ws.Range('E44').Value = '  +4.77%  '
$cell = $ws.Range('D45')
$cell.NumberFormat = "@"
$cell.Value = '0.275'
$cell.Style = "Normal"
$ws.Range('E45').Value = '  +9.60%  '
$ws.Range('E46').Value = '  +8.38%  '
$cell = $ws.Range('D47')
$cell.NumberFormat = "@"
$cell.Value = '124.88'
$cell.Style = "Normal"
$ws.Range('E47').Value = '  +1.28%  '
$cell = $ws.Range('D49')
$cell.NumberFormat = "@"
$cell.Value = '34.63'
$cell.Style = "Normal"
$ws.Range('E50').Value = '  +1.30%  '
$cell = $ws.Range('D51')
$cell.NumberFormat = "@"
$cell.Value = '24.96'
$cell.Style = "Normal"
$ws.Range('E51').Value = '  +6.49%  '
